# This document originally had:
#   P1: "Apartado I " / "Mantenimiento" heading, carrying the _GoBack bookmark
#   P2: empty
#   P3: "En el contrato firmado ... finalización de dicho período. "  (one run)
#   P4: "Durante la fase de mantenimiento se controlará ... son los siguientes:  " (three runs)
#   P5: empty (trailing paragraph)
#
# The edit keeps P1/P2 as-is (minus the bookmark), rewrites P3 into a
# slightly reworded 3-run paragraph, splits P4 into two brand-new
# paragraphs (a fresh "El ciclo de vida ..." paragraph plus a reworded
# "En este proyecto, los tipos de cambios ..." paragraph), and relocates
# the _GoBack bookmark from P1 down onto the trailing empty paragraph
# (P5), which picks up P4's old paragraph formatting.
#
# We use Range.InsertXML on each target paragraph's Range so the exact
# run/paragraph properties from the target markup are produced (far more
# reliable than Find/Replace for this kind of run-splitting/merging).

$d = $word.ActiveDocument

# --- Step 1: move the _GoBack bookmark onto the trailing empty paragraph ---
# (Do this before touching the heading paragraph so the bookmark always
# exists somewhere in the document.)
$p5 = $d.Paragraphs.Item(5)
$p5xml = @'
<w:p><w:pPr><w:spacing w:after="0" w:line="300" w:lineRule="atLeast"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="es-ES"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$p5.Range.InsertXML($p5xml)

# --- Step 2: heading paragraph keeps its two runs, loses the bookmark ---
$p1 = $d.Paragraphs.Item(1)
$p1xml = @'
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="es-ES"/></w:rPr><w:t xml:space="preserve">Apartado I </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="es-ES"/></w:rPr><w:t>Mantenimiento</w:t></w:r></w:p>
'@
$p1.Range.InsertXML($p1xml)

# --- Step 3: reword/re-split the "En el contrato firmado ..." paragraph ---
$p3 = $d.Paragraphs.Item(3)
$p3xml = @'
<w:p><w:pPr><w:spacing w:line="300" w:lineRule="atLeast"/><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:lang w:eastAsia="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:lang w:eastAsia="es-ES"/></w:rPr><w:t>En el contrato firmado con el cliente, además del desarrollo de la aplicación, se ha pactado el mantenimiento de la misma por un período de dos años, prorrogable a partir de la</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:lang w:eastAsia="es-ES"/></w:rPr><w:t xml:space="preserve"> fecha de</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:lang w:eastAsia="es-ES"/></w:rPr><w:t xml:space="preserve"> finalización de dicho período. </w:t></w:r></w:p>
'@
$p3.Range.InsertXML($p3xml)

# --- Step 4: split "Durante la fase ..." into the two new paragraphs ---
$p4 = $d.Paragraphs.Item(4)
$p4xml = @'
<w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>El ciclo de vida del desarrollo del software no termina cuando la aplicación está terminada, instalada y en funcionamiento. Falta la etapa final, que puede llegar a ser la más larga, que es la de mantenimiento. Una vez finalizadas las fases anteriores y con el software funcionando, hay que controlar, mejorar y optimizar</w:t></w:r><w:r><w:t xml:space="preserve"> el software, y </w:t></w:r><w:r><w:t>realizar</w:t></w:r><w:r><w:t xml:space="preserve"> las modificaciones que se consideren necesarias </w:t></w:r><w:r><w:t>para hacer frente a</w:t></w:r><w:r><w:t xml:space="preserve"> las diferentes situaciones que vayan surgiendo</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="300" w:lineRule="atLeast"/><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:lang w:eastAsia="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:lang w:eastAsia="es-ES"/></w:rPr><w:t>En este proyecto, l</w:t></w:r><w:r><w:t xml:space="preserve">os tipos de cambios que se prevén son los siguientes: </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:lang w:eastAsia="es-ES"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$p4.Range.InsertXML($p4xml)

Write-Host "Paragraph count after edit:" $d.Paragraphs.Count
